$d = $word.ActiveDocument

# --- Edit 1: merge the two runs in the second paragraph into one run ---
# The original text is split across two runs: "...work with other " + "sources.' Some..."
# A Find/Replace over the same text naturally collapses it into a single run.
$find1 = "It is understood that public domain cyber-security data in this space is sparse and that you may have to work with other sources.’"
$rng1 = $d.Content
$rng1.Find.ClearFormatting()
$rng1.Find.Execute($find1, $true, $false, $false, $false, $false, $true, 1, $false, $find1, 2) | Out-Null

# --- Edit 2: append new paragraphs describing the Enron dataset decision ---
$insertPoint = $d.Content
$insertPoint.Collapse(0)

$W = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

$newXml = "<w:p $W/>"
$newXml += "<w:p $W><w:pPr><w:ind w:left=`"720`"/></w:pPr><w:r><w:t>Our team decided to go with option 1, the ‘Enron spam e-mails’ given a couple of considerations, such as:</w:t></w:r></w:p>"
$newXml += "<w:p $W><w:r><w:t xml:space=`"preserve`">-         </w:t></w:r><w:r><w:rPr><w:i/><w:iCs/></w:rPr><w:t>The integrity, replicability and continuity of the data</w:t></w:r><w:r><w:t>: The e-mails were formatted in an initial gzip (zipped) state, which was unzipped and resulted in 6 separate folders – determined by the author of the e-mails in each respective folder. These were then placed in a single folder. All e-mails were represented in an individual and uniquely numbered plain text file. A ‘summary’ text file also featured in every of the 6 folders, explaining the correspondence and the separation of spam from ham.</w:t></w:r></w:p>"
$newXml += "<w:p $W><w:r><w:t xml:space=`"preserve`">-        </w:t></w:r><w:r><w:rPr><w:i/><w:iCs/></w:rPr><w:t>Ease of ulterior performance analysis:</w:t></w:r><w:r><w:t xml:space=`"preserve`"> Given the summary files, we possess prior knowledge of what to expect the ratio spam:normal would look like. Moreover, each of the 6 folders had its e-mails separated in 2 subfolders: ‘Spam’ for all spam e-mails and ‘Ham’ for all normal correspondences.</w:t></w:r></w:p>"
$newXml += "<w:p $W><w:r><w:lastRenderedPageBreak/><w:t xml:space=`"preserve`">-        </w:t></w:r><w:r><w:rPr><w:i/><w:iCs/></w:rPr><w:t>Source validity:</w:t></w:r><w:r><w:rPr><w:i/><w:iCs/></w:rPr><w:t xml:space=`"preserve`"> </w:t></w:r><w:r><w:t>The Enron case is well known, and all e-mails were sure to be good sources of insight. Moreover, those already categorized as spam were known to have been deemed so on the basis of something better than an industrial spam filter, e.g: yahoo, gmail etc.</w:t></w:r></w:p>"
$newXml += "<w:p $W><w:r><w:t xml:space=`"preserve`">  The only hardship encountered was uploading the files to GitHub – as even in a compressed state they were slightly (28mb) over the upload limit supported by the browser version (25mb). Our team decided to eliminate one of the 6 folders at random (enron 4, specifically, was arbitrarily chosen) as enough data was available without it as well.</w:t></w:r></w:p>"

$insertPoint.InsertXML($newXml)

Write-Host "Edits applied successfully."
